$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the H3 formula: originally a hardcoded 0, now computed like the rest of the column.
$ws.Range("H3").Formula = "=IF(G3>F3,G3-F3,0)"

# H4:H38 get the corrected (non-off-by-one) formula referencing the same row
# (relative references auto-adjust per row, same as an Excel fill-down).
$ws.Range("H4:H38").Formula = "=IF(G4>F4,G4-F4,0)"

# D36:D38 become individual (unshared) formulas referencing D3, matching D35's pattern.
$ws.Range("D36").Formula = "=D3/2"
$ws.Range("D37").Formula = "=D3/2"
$ws.Range("D38").Formula = "=D3/2"

# Restore the view: scroll back to top-left and move the selection to H40.
$ws.Range("H40").Select()
